# Add a "biling" column to the INSERT statements generated in column H of
# the "domande" sheet and append a trailing "COMMIT;" statement.
#
# Before: INSERT INTO domande(id_domanda, testo, id_categoria, id_immagine, risposta) VALUES (...);
# After:  INSERT INTO domande(id_domanda, testo, id_categoria, id_immagine, risposta, biling) VALUES (...);COMMIT;
#
# This mirrors the CONCAT() formulas stored in H2:H75 (id 1..74), which build
# the SQL insert strings used elsewhere in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("domande")

for ($i = 2; $i -le 75; $i++) {
    $cell = $ws.Range("H$i")
    $oldFormula = $cell.Formula
    if ($oldFormula -ne "") {
        $newFormula = $oldFormula.Replace("risposta) VALUES (", "risposta, biling) VALUES (")
        $newFormula = $newFormula.Replace('");")', '");COMMIT;")')
        $cell.Formula = $newFormula
    }
}

# Match the author's final selection/view state: active cell H1 with the
# whole column H selected.
$ws.Activate()
$ws.Range("H1:H1048576").Select()
